{"js": "// Revert to commit a91223903f182c8607aac0a10b9e7e61a1a8d25d\n//\n// 1) \"...module 'INOUT'\" -> \"...module 'OUT'\" (two occurrences: the\n//    ETPDATA paragraph and the HEADER paragraph). The run(s) spanning\n//    from \"array called\" / \"'HEADER'\" through the end of the paragraph\n//    are collapsed into a single run with the corrected text.\n// 2) The \"desired wave group properties ... rather than read from text\n//    file.\" bullet paragraph is removed entirely.\n// 3) A \"_GoBack\" bookmark (start+end, empty range) is added immediately\n//    before the \"References:\" run.\n\nconst body = context.document.body;\n\n// --- 1a) ETPDATA paragraph: drop the stray \"IN\" so the module name reads 'OUT' ---\nlet results = body.search(\"array called\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const startRange = results.items[0];\n  const para = startRange.paragraphs.getFirst();\n  const paraEnd = para.getRange(\"End\");\n  const fullRange = startRange.expandTo(paraEnd);\n  fullRange.insertText(\"array called \\u2018ETPDATA\\u2019 in the module \\u2018OUT\\u2019\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 1b) HEADER paragraph: same fix ---\nresults = body.search(\"\\u2018HEADER\\u2019\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const startRange = results.items[0];\n  const para = startRange.paragraphs.getFirst();\n  const paraEnd = para.getRange(\"End\");\n  const fullRange = startRange.expandTo(paraEnd);\n  fullRange.insertText(\"\\u2018HEADER\\u2019 in the module \\u2018OUT\\u2019\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) Remove the \"desired wave group properties ...\" bullet paragraph ---\nresults = body.search(\"desired wave group properties\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const para = results.items[0].paragraphs.getFirst();\n  para.delete();\n}\nawait context.sync();\n\n// --- 3) Insert the \"_GoBack\" bookmark right before \"References:\" ---\nresults = body.search(\"References:\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const refPara = results.items[0].paragraphs.getFirst();\n  const startOfPara = refPara.getRange(\"Start\");\n  startOfPara.insertBookmark(\"_GoBack\");\n}\nawait context.sync();\n", "ps1": "# Revert to commit a91223903f182c8607aac0a10b9e7e61a1a8d25d\n#\n# 1) \"...module 'INOUT'\" -> \"...module 'OUT'\" (two occurrences: the\n#    ETPDATA paragraph and the HEADER paragraph) via Find/Replace All,\n#    which naturally collapses the split runs into a single run.\n# 2) The \"desired wave group properties ... rather than read from text\n#    file.\" bullet paragraph is removed entirely.\n# 3) A \"_GoBack\" bookmark (start+end, empty range) is added immediately\n#    before the \"References:\" run.\n\n$d = $word.ActiveDocument\n\n# --- 1) Fix the \"INOUT\" -> \"OUT\" typo in both occurrences ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"INOUT\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"OUT\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# --- 2) Remove the \"desired wave group properties ...\" bullet paragraph ---\n$range = $d.Content\n$find2 = $range.Find\n$find2.Text = \"desired wave group properties\"\n$find2.Execute()\nif ($find2.Found) {\n    $range.Expand(4)  # wdParagraph\n    $range.Delete()\n}\n\n# --- 3) Insert the \"_GoBack\" bookmark right before \"References:\" ---\n$range2 = $d.Content\n$find3 = $range2.Find\n$find3.Text = \"References:\"\n$find3.Execute()\nif ($find3.Found) {\n    $range2.Collapse(1)  # wdCollapseStart\n    $d.Bookmarks.Add(\"_GoBack\", $range2)\n}\n"}
